$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G24").Value = 1.95
$ws.Range("H24").Value = 3.25
$ws.Range("I24").Value = 4.1
$ws.Range("J24").Value = 2.75
$ws.Range("L24").Value = 5
$ws.Range("Q24").Value = 2.6
$ws.Range("R24").Value = 1.48
$ws.Range("U24").Value = 5.5
$ws.Range("V24").Value = 1.14
$ws.Range("AB24").Value = 8
$ws.Range("AD24").Value = 17
$ws.Range("AE24").Value = 21
$ws.Range("AG24").Value = 6.5
$ws.Range("AL24").Value = 8.5
$ws.Range("AM24").Value = 19
$ws.Range("AO24").Value = 41
$ws.Range("AR24").Value = 2
$ws.Range("AS24").Value = 1.85

$ws.Range("G25").Value = 2.1
$ws.Range("H25").Value = 3.25
$ws.Range("I25").Value = 3.6
$ws.Range("J25").Value = 2.88
$ws.Range("K25").Value = 1.95
$ws.Range("L25").Value = 4.5
$ws.Range("N25").Value = 7.5
$ws.Range("S25").Value = 3.8
$ws.Range("Y25").Value = 2.05
$ws.Range("Z25").Value = 1.7
$ws.Range("AA25").Value = 6
$ws.Range("AB25").Value = 9
$ws.Range("AC25").Value = 9.5
$ws.Range("AD25").Value = 19
$ws.Range("AG25").Value = 7.5
$ws.Range("AH25").Value = 6
$ws.Range("AI25").Value = 19
$ws.Range("AJ25").Value = 67
$ws.Range("AK25").Value = 501
$ws.Range("AL25").Value = 8.5
$ws.Range("AM25").Value = 17
$ws.Range("AN25").Value = 13
$ws.Range("AO25").Value = 41
$ws.Range("AP25").Value = 34
$ws.Range("AQ25").Value = 41

$ws.Range("G28").Value = 1.96
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = 1.87
$ws.Range("L28").Value = 5
$ws.Range("M28").Value = 1.13
$ws.Range("N28").Value = 6
$ws.Range("R28").Value = 1.41
$ws.Range("U28").Value = 5.5
$ws.Range("V28").Value = 1.14
$ws.Range("AB28").Value = 8
$ws.Range("AJ28").Value = 81
$ws.Range("AO28").Value = 41

$ws.Range("O30").Value = 1.73
$ws.Range("P30").Value = 2
$ws.Range("Q30").Value = 3.5
$ws.Range("R30").Value = 1.3
$ws.Range("U30").Value = 8
$ws.Range("V30").Value = 1.08
$ws.Range("W30").Value = 1.75
$ws.Range("X30").Value = 2.05

$ws.Range("G31").Value = 2.25
$ws.Range("H31").Value = 2.9
$ws.Range("I31").Value = 3.75
$ws.Range("J31").Value = 3.2
$ws.Range("K31").Value = 1.8
$ws.Range("O31").Value = 1.67
$ws.Range("P31").Value = 2.1
$ws.Range("AB31").Value = 9
$ws.Range("AD31").Value = 21
$ws.Range("AE31").Value = 26
$ws.Range("AG31").Value = 5
$ws.Range("AH31").Value = 6

$ws.Range("G50").Value = 2.25
$ws.Range("K50").Value = 1.95
$ws.Range("M50").Value = 1.11
$ws.Range("N50").Value = 6.5
$ws.Range("Q50").Value = 2.6
$ws.Range("R50").Value = 1.48
$ws.Range("AR50").Value = 1.98
$ws.Range("AS50").Value = 1.88

$ws.Range("G51").Value = 1.44
$ws.Range("H51").Value = 4.2
$ws.Range("J51").Value = 2
$ws.Range("R51").Value = 1.75
$ws.Range("AD51").Value = 9
$ws.Range("AH51").Value = 8.5
$ws.Range("AL51").Value = 15
$ws.Range("AO51").Value = 81
$ws.Range("AP51").Value = 51
$ws.Range("AS51").Value = 2.43

$ws.Range("G84").Value = 1.5
$ws.Range("H84").Value = 4
$ws.Range("I84").Value = 7.5
$ws.Range("L84").Value = 7
$ws.Range("Y84").Value = 2.2
$ws.Range("Z84").Value = 1.62
$ws.Range("AA84").Value = 5.5
$ws.Range("AD84").Value = 10
$ws.Range("AF84").Value = 34
$ws.Range("AG84").Value = 9
$ws.Range("AH84").Value = 8
$ws.Range("AJ84").Value = 81
$ws.Range("AM84").Value = 34
$ws.Range("AN84").Value = 21
$ws.Range("AO84").Value = 81

$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 3.5
$ws.Range("I94").Value = 3.5
$ws.Range("J94").Value = 2.75
$ws.Range("N94").Value = 9.5
$ws.Range("O94").Value = 1.3
$ws.Range("P94").Value = 3.4
$ws.Range("R94").Value = 1.8
$ws.Range("AA94").Value = 7
$ws.Range("AB94").Value = 9.5
$ws.Range("AD94").Value = 17
$ws.Range("AE94").Value = 17
$ws.Range("AG94").Value = 9.5
$ws.Range("AL94").Value = 10
$ws.Range("AO94").Value = 41
$ws.Range("AQ94").Value = 41
